$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1325
$ws.Range("H69").Value = 3333
$ws.Range("I69").Value = 3333
$ws.Range("K69").Value = 9999
$ws.Range("M69").Value = -9125
$ws.Range("H72").Value = 3333
$ws.Range("I72").Value = 3333
$ws.Range("K72").Value = 29997
$ws.Range("M72").Value = -25629
$ws.Range("H80").Value = 1200
$ws.Range("J80").Value = 1200
$ws.Range("L80").Value = 3600
$ws.Range("N80").Value = -5596
$ws.Range("H83").Value = 1200
$ws.Range("J83").Value = 1200
$ws.Range("L83").Value = 10800
$ws.Range("N83").Value = -20784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2809.8572
$ws.Range("I61").Value = 2999.8
$ws.Range("J61").Value = 2335
$ws.Range("K61").Value = 2999.8
$ws.Range("L61").Value = 2335
$ws.Range("M61").Value = -2787.8
$ws.Range("N61").Value = -2759
$ws.Range("H88").Value = 2336.1667
$ws.Range("I88").Value = 1373.3334
$ws.Range("K88").Value = 1373.3334
$ws.Range("M88").Value = -967.3334
$ws.Range("H91").Value = 2336.1667
$ws.Range("I91").Value = 1373.3334
$ws.Range("K91").Value = 1373.3334
$ws.Range("M91").Value = 30.66660000000002
$ws.Range("H136").Value = 2809.8572
$ws.Range("I136").Value = 2999.8
$ws.Range("J136").Value = 2335
$ws.Range("K136").Value = 8999.400000000001
$ws.Range("L136").Value = 7005
$ws.Range("M136").Value = -6449.400000000001
$ws.Range("N136").Value = -12105

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86622
$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -263112
$ws.Range("H82").Value = 18375.834
$ws.Range("I82").Value = 14051
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 14051
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = -13668
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 18375.834
$ws.Range("I85").Value = 14051
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 14051
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = -12725
$ws.Range("N85").Value = -42652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 6155.3076
$ws.Range("I33").Value = 1018.875
$ws.Range("J33").Value = 14373.6
$ws.Range("K33").Value = 1018.875
$ws.Range("L33").Value = 14373.6
$ws.Range("M33").Value = -639.875
$ws.Range("N33").Value = -15131.6
$ws.Range("H44").Value = 29999.75
$ws.Range("I44").Value = 30000
$ws.Range("J44").Value = 29999.666
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 29999.666
$ws.Range("M44").Value = -29558
$ws.Range("N44").Value = -30883.666
$ws.Range("H51").Value = 32001.5
$ws.Range("J51").Value = 27668.666
$ws.Range("L51").Value = 27668.666
$ws.Range("N51").Value = -29140.666
$ws.Range("H60").Value = 29331.25
$ws.Range("J60").Value = 34744.332
$ws.Range("L60").Value = 34744.332
$ws.Range("N60").Value = -35766.332
$ws.Range("H61").Value = 32001.5
$ws.Range("J61").Value = 27668.666
$ws.Range("L61").Value = 27668.666
$ws.Range("N61").Value = -28364.666
$ws.Range("H74").Value = 87543
$ws.Range("J74").Value = 87543
$ws.Range("L74").Value = 87543
$ws.Range("N74").Value = -89291
$ws.Range("H77").Value = 87543
$ws.Range("J77").Value = 87543
$ws.Range("L77").Value = 262629
$ws.Range("N77").Value = -271365
$ws.Range("H99").Value = 1432591.9
$ws.Range("I99").Value = 1251285.8
$ws.Range("K99").Value = 1251285.8
$ws.Range("M99").Value = -1249787.8
$ws.Range("H126").Value = 1432591.9
$ws.Range("I126").Value = 1251285.8
$ws.Range("K126").Value = 3753857.4
$ws.Range("M126").Value = -3751387.4
$ws.Range("H132").Value = 1468.5
$ws.Range("I132").Value = 1362.2
$ws.Range("K132").Value = 4086.6
$ws.Range("M132").Value = -1556.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 110.454544
$ws.Range("I7").Value = 109.2
$ws.Range("J7").Value = 111.5
$ws.Range("K7").Value = 327.6
$ws.Range("L7").Value = 334.5
$ws.Range("M7").Value = -215.6
$ws.Range("N7").Value = -558.5
$ws.Range("H17").Value = 350.25
$ws.Range("J17").Value = 466.66666
$ws.Range("L17").Value = 1399.99998
$ws.Range("N17").Value = -1737.99998
$ws.Range("H35").Value = 11105
$ws.Range("I35").Value = 457.5
$ws.Range("K35").Value = 1372.5
$ws.Range("M35").Value = -1084.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 108.833336
$ws.Range("I2").Value = 126.1
$ws.Range("K2").Value = 126.1
$ws.Range("M2").Value = -13.09999999999999
$ws.Range("H107").Value = 1383.3334
$ws.Range("I107").Value = 1304.8334
$ws.Range("J107").Value = 1540.3334
$ws.Range("K107").Value = 1304.8334
$ws.Range("L107").Value = 1540.3334
$ws.Range("M107").Value = 615.1666
$ws.Range("N107").Value = -5380.3334
$ws.Range("H119").Value = 78000
$ws.Range("J119").Value = 78000
$ws.Range("L119").Value = 78000
$ws.Range("N119").Value = -87676
$ws.Range("H123").Value = 89497
$ws.Range("J123").Value = 89497
$ws.Range("L123").Value = 89497
$ws.Range("N123").Value = -94397

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("H26").Value = 3562.125
$ws.Range("I26").Value = 1165.6666
$ws.Range("K26").Value = 1165.6666
$ws.Range("M26").Value = -870.6666
$ws.Range("H31").Value = 3627.3076
$ws.Range("I31").Value = 1213.75
$ws.Range("J31").Value = 4700
$ws.Range("K31").Value = 1213.75
$ws.Range("L31").Value = 4700
$ws.Range("M31").Value = -965.75
$ws.Range("N31").Value = -5196
$ws.Range("H55").Value = 411.33334
$ws.Range("I55").Value = 422.75
$ws.Range("J55").Value = 407.18182
$ws.Range("K55").Value = 422.75
$ws.Range("L55").Value = 407.18182
$ws.Range("M55").Value = -249.75
$ws.Range("N55").Value = -753.18182
$ws.Range("H76").Value = 15187
$ws.Range("J76").Value = 15187
$ws.Range("L76").Value = 15187
$ws.Range("N76").Value = -15863
$ws.Range("H79").Value = 15187
$ws.Range("J79").Value = 15187
$ws.Range("L79").Value = 15187
$ws.Range("N79").Value = -17527
$ws.Range("H100").Value = 4984.5
$ws.Range("J100").Value = 4970
$ws.Range("L100").Value = 4970
$ws.Range("N100").Value = -6052
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
